$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.779.85'
$ws.Range("E2").Value = '  -0.52%  '
$ws.Range("D3").Value = '3.850.14'
$ws.Range("E3").Value = '  +2.16%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '423.91'
$ws.Range("E5").Value = '  +0.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.21'
$ws.Range("E6").Value = '  -1.31%  '
$ws.Range("D7").Value = '3.839.77'
$ws.Range("E7").Value = '  +2.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.608'
$ws.Range("E8").Value = '  -6.20%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.724'
$ws.Range("E10").Value = '  -6.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.165'
$ws.Range("E11").Value = '  -9.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000350'
$ws.Range("E12").Value = '  -15.23%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.62'
$ws.Range("E13").Value = '  -5.47%  '
$ws.Range("B14").Value = 'Uniswap'
$ws.Range("C14").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.57'
$ws.Range("E14").Value = '  +24.98%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '4.474.41'
$ws.Range("E15").Value = '  +2.66%  '
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '10.00'
$ws.Range("E16").Value = '  -3.78%  '
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.138'
$ws.Range("E17").Value = '  -1.40%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.831.68'
$ws.Range("E18").Value = '  +1.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '19.64'
$ws.Range("E19").Value = '  -5.19%  '
$ws.Range("D20").Value = '67.011.28'
$ws.Range("E20").Value = '  -0.20%  '
$ws.Range("E21").Value = '  -6.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '405.03'
$ws.Range("E22").Value = '  -9.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.43'
$ws.Range("E23").Value = '  -9.85%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.63'
$ws.Range("E24").Value = '  -5.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.02'
$ws.Range("E25").Value = '  -3.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '37.16'
$ws.Range("E26").Value = '  -4.66%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.70'
$ws.Range("E27").Value = '  +11.59%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.25'
$ws.Range("E28").Value = '  -2.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.52'
$ws.Range("E29").Value = '  -6.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '693.44'
$ws.Range("E30").Value = '  +2.79%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.77'
$ws.Range("E31").Value = '  +2.09%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.121'
$ws.Range("E32").Value = '  -3.99%  '
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.34'
$ws.Range("E33").Value = '  -2.80%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.31'
$ws.Range("E34").Value = '  +0.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.152'
$ws.Range("E35").Value = '  -8.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '38.29'
$ws.Range("E36").Value = '  -8.78%  '
$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.13'
$ws.Range("E38").Value = '  -3.16%  '
$ws.Range("B39").Value = 'PEPE'
$ws.Range("C39").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D39").Value = '0.0₃0783'
$ws.Range("E39").Value = '  +1.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0457'
$ws.Range("E40").Value = '  -7.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.95'
$ws.Range("E41").Value = '  -0.95%  '
$ws.Range("E42").Value = '  +0.49%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.135'
$ws.Range("E43").Value = '  -9.56%  '
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.45'
$ws.Range("E44").Value = '  +2.13%  '
$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '146.05'
$ws.Range("E45").Value = '  -1.08%  '
$ws.Range("B46").Value = 'LidoDAOToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.28'
$ws.Range("E46").Value = '  -5.29%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.55'
$ws.Range("E47").Value = '  -7.67%  '
$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.09'
$ws.Range("E48").Value = '  -2.70%  '
$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.06'
$ws.Range("E49").Value = '  -5.25%  '
$ws.Range("B50").Value = 'WEMIXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.56'
$ws.Range("E50").Value = '  -3.63%  '
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.78'
$ws.Range("E51").Value = '  -4.11%  '
